$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# --- Row 10: height tweak ---
$ws.Rows.Item(10).RowHeight = 20.5

# --- Row 16: fix_prices_to_2030 -> fix_fuel_prices_to_year, now enabled, updated description ---
$ws.Range("A16").Value = "fix_fuel_prices_to_year"
$ws.Range("B16").Value = $true
$ws.Range("C16").Value = "for verification runs. Fix fuel prices, CO2 and electricity demand to "

# --- Row 17: fix_prices_to_2020 -> fix_price_year, now a year number instead of a boolean ---
$ws.Range("A17").Value = "fix_price_year"
$ws.Range("B17").Value = 2020
$ws.Range("C17").Value = "to this year"

# --- Row 19: clarify description ---
$ws.Range("C19").Value = "so far this is only for NL. If False"

# --- Row 26 (install_at_look_ahead_year) removed entirely ---
$ws.Rows.Item(26).Delete()

# --- New row 28: formula warning when prices are fixed with fuel trends also on ---
$ws.Rows.Item(28).Insert()
$ws.Range("B28").Formula = '=IF(AND(B16=TRUE,B13>0),"PRICES are fixed, no fuel trends are considered","ok")'

# --- New row 29: AMIRIS-changed sanity-check formula (references now-removed cell -> #REF!) ---
$ws.Rows.Item(29).Insert()
$ws.Range("B29").Formula = '=IF(#REF!=TRUE,"DANGER!!!!!","ok")'
$ws.Range("C29").Value = "AMIRIS has changed"

# Remove the old yellow highlight style from the DANGER-check rows (now B30:B33)
$ws.Range("B30:B33").ClearFormats()

# --- Conditional formatting: add a red-fill rule over the new formula block ---
$newRule = $ws.Range("B28:B33").FormatConditions.Add(1, 4, """ok""")
$newRule.Interior.Color = 255
$newRule.SetFirstPriority()

# --- View bookkeeping to mirror the saved workbook state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C20").Select()
